$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data region (rows 2-28) before rewriting with the reorganised / expanded layout
$ws.Range("A2:D42").ClearContents()

$ws.Cells.Item(2,1).Value = 45
$ws.Cells.Item(2,2).Value = "water bottle"
$ws.Cells.Item(2,3).Value = 240000
$ws.Cells.Item(2,4).Value = "Group 1"

$ws.Cells.Item(3,1).Value = 523
$ws.Cells.Item(3,2).Value = "hoodies"
$ws.Cells.Item(3,3).Value = 80000
$ws.Cells.Item(3,4).Value = "Group 1"

$ws.Cells.Item(4,3).Value = 320000
$ws.Cells.Item(4,4).Value = "Total Group 1"

$ws.Cells.Item(5,1).Value = 2
$ws.Cells.Item(5,2).Value = "bags"
$ws.Cells.Item(5,3).Value = 240000
$ws.Cells.Item(5,4).Value = "Group 2"

$ws.Cells.Item(6,1).Value = 23542
$ws.Cells.Item(6,2).Value = "wings"
$ws.Cells.Item(6,3).Value = 80000
$ws.Cells.Item(6,4).Value = "Group 2"

$ws.Cells.Item(7,3).Value = 320000
$ws.Cells.Item(7,4).Value = "Total Group 2"

$ws.Cells.Item(8,1).Value = 4567
$ws.Cells.Item(8,2).Value = "suspension"
$ws.Cells.Item(8,3).Value = 200000
$ws.Cells.Item(8,4).Value = "Group 3"

$ws.Cells.Item(9,1).Value = 45
$ws.Cells.Item(9,2).Value = "rice bags"
$ws.Cells.Item(9,3).Value = 120000
$ws.Cells.Item(9,4).Value = "Group 3"

$ws.Cells.Item(10,3).Value = 320000
$ws.Cells.Item(10,4).Value = "Total Group 3"

$ws.Cells.Item(11,1).Value = 123123
$ws.Cells.Item(11,2).Value = "Dog"
$ws.Cells.Item(11,3).Value = 120000
$ws.Cells.Item(11,4).Value = "Group 4"

$ws.Cells.Item(12,1).Value = 3457
$ws.Cells.Item(12,2).Value = "snake"
$ws.Cells.Item(12,3).Value = 90000
$ws.Cells.Item(12,4).Value = "Group 4"

$ws.Cells.Item(13,1).Value = 23452
$ws.Cells.Item(13,2).Value = "Screens"
$ws.Cells.Item(13,3).Value = 80000
$ws.Cells.Item(13,4).Value = "Group 4"

$ws.Cells.Item(14,1).Value = 5
$ws.Cells.Item(14,2).Value = "bottle"
$ws.Cells.Item(14,3).Value = 30000
$ws.Cells.Item(14,4).Value = "Group 4"

$ws.Cells.Item(15,3).Value = 320000
$ws.Cells.Item(15,4).Value = "Total Group 4"

$ws.Cells.Item(16,1).Value = 567
$ws.Cells.Item(16,2).Value = "boots"
$ws.Cells.Item(16,3).Value = 80000
$ws.Cells.Item(16,4).Value = "Group 5"

$ws.Cells.Item(17,1).Value = 345
$ws.Cells.Item(17,2).Value = "Weights"
$ws.Cells.Item(17,3).Value = 65000
$ws.Cells.Item(17,4).Value = "Group 5"

$ws.Cells.Item(18,1).Value = 74
$ws.Cells.Item(18,2).Value = "watches"
$ws.Cells.Item(18,3).Value = 60000
$ws.Cells.Item(18,4).Value = "Group 5"

$ws.Cells.Item(19,1).Value = 234
$ws.Cells.Item(19,2).Value = "mouse"
$ws.Cells.Item(19,3).Value = 60000
$ws.Cells.Item(19,4).Value = "Group 5"

$ws.Cells.Item(20,1).Value = 24
$ws.Cells.Item(20,2).Value = "Grape"
$ws.Cells.Item(20,3).Value = 40000
$ws.Cells.Item(20,4).Value = "Group 5"

$ws.Cells.Item(21,1).Value = 4363
$ws.Cells.Item(21,2).Value = "rags"
$ws.Cells.Item(21,3).Value = 15000
$ws.Cells.Item(21,4).Value = "Group 5"

$ws.Cells.Item(22,3).Value = 320000
$ws.Cells.Item(22,4).Value = "Total Group 5"

$ws.Cells.Item(23,1).Value = 2131
$ws.Cells.Item(23,2).Value = "Cat"
$ws.Cells.Item(23,3).Value = 60000
$ws.Cells.Item(23,4).Value = "Group 6"

$ws.Cells.Item(24,1).Value = 576
$ws.Cells.Item(24,2).Value = "hippo"
$ws.Cells.Item(24,3).Value = 60000
$ws.Cells.Item(24,4).Value = "Group 6"

$ws.Cells.Item(25,1).Value = 245
$ws.Cells.Item(25,2).Value = "tea bags"
$ws.Cells.Item(25,3).Value = 60000
$ws.Cells.Item(25,4).Value = "Group 6"

$ws.Cells.Item(26,1).Value = 2435345
$ws.Cells.Item(26,2).Value = "spoons"
$ws.Cells.Item(26,3).Value = 40000
$ws.Cells.Item(26,4).Value = "Group 6"

$ws.Cells.Item(27,1).Value = 4363
$ws.Cells.Item(27,2).Value = "rags"
$ws.Cells.Item(27,3).Value = 100000
$ws.Cells.Item(27,4).Value = "Group 6"

$ws.Cells.Item(28,3).Value = 320000
$ws.Cells.Item(28,4).Value = "Total Group 6"

$ws.Cells.Item(29,1).Value = 4363
$ws.Cells.Item(29,2).Value = "rags"
$ws.Cells.Item(29,3).Value = 320000
$ws.Cells.Item(29,4).Value = "Group 7"

$ws.Cells.Item(30,3).Value = 320000
$ws.Cells.Item(30,4).Value = "Total Group 7"

$ws.Cells.Item(31,1).Value = 456
$ws.Cells.Item(31,2).Value = "yoyo"
$ws.Cells.Item(31,3).Value = 320000
$ws.Cells.Item(31,4).Value = "Group 8"

$ws.Cells.Item(32,3).Value = 320000
$ws.Cells.Item(32,4).Value = "Total Group 8"

$ws.Cells.Item(33,1).Value = 4363
$ws.Cells.Item(33,2).Value = "rags"
$ws.Cells.Item(33,3).Value = 320000
$ws.Cells.Item(33,4).Value = "Group 9"

$ws.Cells.Item(34,3).Value = 320000
$ws.Cells.Item(34,4).Value = "Total Group 9"

$ws.Cells.Item(35,1).Value = 4363
$ws.Cells.Item(35,2).Value = "rags"
$ws.Cells.Item(35,3).Value = 95000
$ws.Cells.Item(35,4).Value = "Group 10"

$ws.Cells.Item(36,1).Value = 56
$ws.Cells.Item(36,2).Value = "hats"
$ws.Cells.Item(36,3).Value = 225000
$ws.Cells.Item(36,4).Value = "Group 10"

$ws.Cells.Item(37,3).Value = 320000
$ws.Cells.Item(37,4).Value = "Total Group 10"

$ws.Cells.Item(38,1).Value = 56
$ws.Cells.Item(38,2).Value = "hats"
$ws.Cells.Item(38,3).Value = 135000
$ws.Cells.Item(38,4).Value = "Group 11"

$ws.Cells.Item(39,1).Value = 456
$ws.Cells.Item(39,2).Value = "yoyo"
$ws.Cells.Item(39,3).Value = 185000
$ws.Cells.Item(39,4).Value = "Group 11"

$ws.Cells.Item(40,3).Value = 320000
$ws.Cells.Item(40,4).Value = "Total Group 11"

$ws.Cells.Item(41,1).Value = 456
$ws.Cells.Item(41,2).Value = "yoyo"
$ws.Cells.Item(41,3).Value = 145000
$ws.Cells.Item(41,4).Value = "Group 12"

$ws.Cells.Item(42,3).Value = 145000
$ws.Cells.Item(42,4).Value = "Total Group 12"

# Apply the highlighted "Total" style (solid yellow fill) to the subtotal rows
$ws.Range("A4:D4").Interior.Color = 65535
$ws.Range("A7:D7").Interior.Color = 65535
$ws.Range("A10:D10").Interior.Color = 65535
$ws.Range("A15:D15").Interior.Color = 65535
$ws.Range("A22:D22").Interior.Color = 65535
$ws.Range("A28:D28").Interior.Color = 65535
$ws.Range("A30:D30").Interior.Color = 65535
$ws.Range("A32:D32").Interior.Color = 65535
$ws.Range("A34:D34").Interior.Color = 65535
$ws.Range("A37:D37").Interior.Color = 65535
$ws.Range("A40:D40").Interior.Color = 65535
$ws.Range("A42:D42").Interior.Color = 65535

"done"
